$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 29004.25
$ws.Range("I46").Value = 1017
$ws.Range("J46").Value = 38333.332
$ws.Range("K46").Value = 3051
$ws.Range("L46").Value = 114999.996
$ws.Range("M46").Value = -2932
$ws.Range("N46").Value = -115237.996

$ws.Range("H60").Value = 29004.25
$ws.Range("I60").Value = 1017
$ws.Range("J60").Value = 38333.332
$ws.Range("K60").Value = 3051
$ws.Range("L60").Value = 114999.996
$ws.Range("M60").Value = -2567
$ws.Range("N60").Value = -115967.996

$ws.Range("H64").Value = 3038.5356
$ws.Range("I64").Value = 2818.75
$ws.Range("J64").Value = 3126.45
$ws.Range("K64").Value = 2818.75
$ws.Range("L64").Value = 3126.45
$ws.Range("M64").Value = -2570.75
$ws.Range("N64").Value = -3622.45

$ws.Range("H67").Value = 3038.5356
$ws.Range("I67").Value = 2818.75
$ws.Range("J67").Value = 3126.45
$ws.Range("K67").Value = 2818.75
$ws.Range("L67").Value = 3126.45
$ws.Range("M67").Value = -1960.75
$ws.Range("N67").Value = -4842.45

$ws.Range("H98").Value = 1237.3513
$ws.Range("I98").Value = 845.26666
$ws.Range("J98").Value = 2917.7144
$ws.Range("K98").Value = 845.26666
$ws.Range("L98").Value = 2917.7144
$ws.Range("M98").Value = 652.73334
$ws.Range("N98").Value = -5913.7144

$ws.Range("H122").Value = 1237.3513
$ws.Range("I122").Value = 845.26666
$ws.Range("J122").Value = 2917.7144
$ws.Range("K122").Value = 2535.79998
$ws.Range("L122").Value = 8753.143199999999
$ws.Range("M122").Value = -85.79997999999978
$ws.Range("N122").Value = -13653.1432

$ws.Range("H135").Value = 47620084
$ws.Range("I135").Value = 785.8182
$ws.Range("J135").Value = 100001310
$ws.Range("K135").Value = 7072.3638
$ws.Range("L135").Value = 900011790
$ws.Range("M135").Value = -4537.3638
$ws.Range("N135").Value = -900016860

$ws.Range("H137").Value = 1588743.2
$ws.Range("I137").Value = 1599.6364
$ws.Range("J137").Value = 4274678.5
$ws.Range("K137").Value = 4798.9092
$ws.Range("L137").Value = 12824035.5
$ws.Range("M137").Value = -2248.9092
$ws.Range("N137").Value = -12829135.5

$ws.Range("H138").Value = 3392842
$ws.Range("I138").Value = 1648.409
$ws.Range("J138").Value = 5409227.5
$ws.Range("K138").Value = 4945.227000000001
$ws.Range("L138").Value = 16227682.5
$ws.Range("M138").Value = 194.7729999999992
$ws.Range("N138").Value = -16237962.5

$ws.Range("H141").Value = 899.4857
$ws.Range("I141").Value = 778.8823
$ws.Range("K141").Value = 2336.6469
$ws.Range("M141").Value = 2843.3531

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 10400
$ws.Range("I34").Value = 9500
$ws.Range("K34").Value = 9500
$ws.Range("M34").Value = -9229

$ws.Range("H74").Value = 20216.154
$ws.Range("I74").Value = 25339.586
$ws.Range("J74").Value = 7859.647
$ws.Range("K74").Value = 25339.586
$ws.Range("L74").Value = 7859.647
$ws.Range("M74").Value = -24465.586
$ws.Range("N74").Value = -9607.647000000001

$ws.Range("H77").Value = 20216.154
$ws.Range("I77").Value = 25339.586
$ws.Range("J77").Value = 7859.647
$ws.Range("K77").Value = 126697.93
$ws.Range("L77").Value = 39298.235
$ws.Range("M77").Value = -122329.93
$ws.Range("N77").Value = -48034.235

$ws.Range("H132").Value = 1907.5264
$ws.Range("I132").Value = 1453.0714
$ws.Range("J132").Value = 3180
$ws.Range("K132").Value = 4359.2142
$ws.Range("L132").Value = 9540
$ws.Range("M132").Value = -1829.2142
$ws.Range("N132").Value = -14600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2324.318
$ws.Range("I105").Value = 2168.261
$ws.Range("J105").Value = 2495.238
$ws.Range("K105").Value = 2168.261
$ws.Range("L105").Value = 2495.238
$ws.Range("M105").Value = -421.261
$ws.Range("N105").Value = -5989.237999999999

$ws.Range("H134").Value = 659523.3
$ws.Range("I134").Value = 1029429.3
$ws.Range("J134").Value = 3780.8635
$ws.Range("K134").Value = 3088287.9
$ws.Range("L134").Value = 11342.5905
$ws.Range("M134").Value = -3085752.9
$ws.Range("N134").Value = -16412.5905

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 851.3
$ws.Range("I16").Value = 618.8333
$ws.Range("J16").Value = 1200
$ws.Range("K16").Value = 618.8333
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = -331.8333
$ws.Range("N16").Value = -1774

$ws.Range("H33").Value = 6269
$ws.Range("I33").Value = 1021.7143
$ws.Range("J33").Value = 43000
$ws.Range("K33").Value = 1021.7143
$ws.Range("L33").Value = 43000
$ws.Range("M33").Value = -642.7143
$ws.Range("N33").Value = -43758

$ws.Range("H58").Value = 6785.421
$ws.Range("I58").Value = 8608.786
$ws.Range("J58").Value = 1680
$ws.Range("K58").Value = 8608.786
$ws.Range("L58").Value = 1680
$ws.Range("M58").Value = -8405.786
$ws.Range("N58").Value = -2086

$ws.Range("H113").Value = 851.3
$ws.Range("I113").Value = 618.8333
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 618.8333
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 1551.1667
$ws.Range("N113").Value = -5540

$ws.Range("H136").Value = 6785.421
$ws.Range("I136").Value = 8608.786
$ws.Range("J136").Value = 1680
$ws.Range("K136").Value = 25826.358
$ws.Range("L136").Value = 5040
$ws.Range("M136").Value = -23276.358
$ws.Range("N136").Value = -10140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 800.13336
$ws.Range("J34").Value = 960
$ws.Range("L34").Value = 2880
$ws.Range("N34").Value = -3048

$ws.Range("H131").Value = 923.59
$ws.Range("J131").Value = 923.59
$ws.Range("L131").Value = 2770.77
$ws.Range("N131").Value = -12850.77

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 734.0625
$ws.Range("I22").Value = 350.5
$ws.Range("J22").Value = 788.8570999999999
$ws.Range("K22").Value = 350.5
$ws.Range("L22").Value = 788.8570999999999
$ws.Range("M22").Value = -55.5
$ws.Range("N22").Value = -1378.8571

$ws.Range("H27").Value = 734.0625
$ws.Range("I27").Value = 350.5
$ws.Range("J27").Value = 788.8570999999999
$ws.Range("K27").Value = 350.5
$ws.Range("L27").Value = 788.8570999999999
$ws.Range("M27").Value = -243.5
$ws.Range("N27").Value = -1002.8571

$ws.Range("H132").Value = 4356.619
$ws.Range("I132").Value = 4632.2
$ws.Range("J132").Value = 3951.353
$ws.Range("K132").Value = 13896.6
$ws.Range("L132").Value = 11854.059
$ws.Range("M132").Value = -11366.6
$ws.Range("N132").Value = -16914.059

$ws.Range("H136").Value = 1908.5385
$ws.Range("I136").Value = 1334.2916
$ws.Range("J136").Value = 2827.3333
$ws.Range("K136").Value = 4002.8748
$ws.Range("L136").Value = 8481.999899999999
$ws.Range("M136").Value = -1452.8748
$ws.Range("N136").Value = -13581.9999
